$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.666.10"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.596.76"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'211.42"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.246"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").Value = "'19.45"
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").Value = "'0.0841"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "1.821.02"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "1.648.63"
$ws.Range("E13").Value = "  +3.19%  "
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "'0.522"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "'65.02"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "26.639.26"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "0.0₃0737"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "'208.91"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  +4.64%  "
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("E23").Value = "  +2.14%  "
$ws.Range("D24").Value = "'8.98"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "'7.11"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").Value = "'15.29"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").Value = "1.288.32"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("E35").Value = "  -5.60%  "
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").Value = "'1.48"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").Value = "'0.829"
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("D40").Value = "'1.03"
$ws.Range("E40").Value = "  +16.03%  "
$ws.Range("D41").Value = "'5.45"
$ws.Range("E41").Value = "  +0.94%  "
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").Value = "'0.782"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("E44").Value = "  -0.93%  "
$ws.Range("D45").Value = "1.732.48"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").Value = "'91.14"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("E47").Value = "  -2.95%  "
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "'7.36"
$ws.Range("E51").Value = "  -1.53%  "
